$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.615.72"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.35%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.849.11"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.15%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.006"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.43%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "314.55"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.06%  "

# Row 6
$ws.Range("E6").Value = "  +0.09%  "

# Row 7
$ws.Range("E7").Value = "  -2.24%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3652"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.57%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "44.70"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.64%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07321"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.22%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.8827"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -5.74%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "20.73"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.85%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.864.59"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.92%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.349"
$ws.Range("D14").Style = "Normal"

# Row 15
$ws.Range("E15").Value = "  -2.91%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.06911"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.53%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.005"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.17%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "78.89"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -3.11%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000008895"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.55%  "

# Row 20
$ws.Range("E20").Value = "  +0.10%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "15.43"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.86%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "27.618.82"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.31%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.994"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.53%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "10.64"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -3.51%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.072.32"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -4.16%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.968"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.31%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "153.67"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.19%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.03"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.45%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "122.01"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +7.41%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.261"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -6.57%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.919"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +12.13%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.08950"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.97%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.7625"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -6.34%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.578"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -5.01%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.977"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.38%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.102"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -6.98%  "

# Row 37
$ws.Range("E37").Value = "  -0.03%  "

# Row 38
$ws.Range("B38").Value = "Hedera"
$ws.Range("C38").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.05378"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.58%  "

# Row 39
$ws.Range("B39").Value = "TrustWalletToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.096"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.71%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01951"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.40%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.807"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -5.22%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.927"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.31%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.5107"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.12%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.1654"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.79%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.276"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -5.72%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.06575"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.67%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.4762"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.63%  "

# Row 48
$ws.Range("E48").Value = "  -1.87%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "104.54"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.14%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.002"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.09%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.630"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.81%  "
